$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12 (I12:L12) - moment analysis PH values
$ws.Range("I12").Value = 0.368443878262713
$ws.Range("J12").Value = 0.03920614569937778
$ws.Range("K12").Value = -1.127910942023814
$ws.Range("L12").Value = 3.329723141697085

# Row 14 (I14:L14) - moment analysis PH values
$ws.Range("I14").Value = 0.3802977549935186
$ws.Range("J14").Value = 0.03528897824179948
$ws.Range("K14").Value = -1.329203199933692
$ws.Range("L14").Value = 3.915931171864705
